$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '63.459.74'
$ws.Range("E2").Value = '  +1.15%  '
# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.667.89'
$ws.Range("E3").Value = '  +3.71%  '
# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.998'
$ws.Range("E4").Value = '  -0.15%  '
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '594.64'
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '143.95'
$ws.Range("E6").Value = '  +0.20%  '
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.998'
$ws.Range("E7").Value = '  -0.16%  '
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.588'
$ws.Range("E8").Value = '  -0.01%  '
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.667.83'
$ws.Range("E9").Value = '  +3.73%  '
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.107'
$ws.Range("E10").Value = '  +0.95%  '
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.69'
$ws.Range("E11").Value = '  +2.56%  '
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.153'
$ws.Range("E12").Value = '  +0.86%  '
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.356'
$ws.Range("E13").Value = '  +2.03%  '
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.53'
$ws.Range("E14").Value = '  +2.70%  '
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.140.08'
$ws.Range("E15").Value = '  +3.54%  '
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '63.331.65'
$ws.Range("E16").Value = '  +1.03%  '
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000145'
$ws.Range("E17").Value = '  +0.42%  '
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.664.13'
$ws.Range("E18").Value = '  +3.80%  '
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.54'
$ws.Range("E19").Value = '  +4.27%  '
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.42'
$ws.Range("E20").Value = '  +2.23%  '
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '339.67'
$ws.Range("E21").Value = '  -0.01%  '
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.82'
$ws.Range("E22").Value = '  +2.81%  '
# Row 23
$ws.Range("E23").Value = '  +0.01%  '
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '67.44'
$ws.Range("E24").Value = '  +0.84%  '
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.69'
$ws.Range("E25").Value = '  +7.17%  '
# Row 26
$ws.Range("B26").Value = 'SuiNetwork'
$ws.Range("C26").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.54'
$ws.Range("E26").Value = '  +2.04%  '
# Row 27
$ws.Range("B27").Value = 'Kaspa'
$ws.Range("C27").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.166'
$ws.Range("E27").Value = '  +1.54%  '
# Row 28
$ws.Range("B28").Value = 'InternetComputer(DFINITY)'
$ws.Range("C28").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.49'
$ws.Range("E28").Value = '  +3.61%  '
# Row 29
$ws.Range("B29").Value = 'Binance-PegBSC-USD'
$ws.Range("C29").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.00'
$ws.Range("E29").Value = '  +0.70%  '
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '539.83'
$ws.Range("E30").Value = '  +18.82%  '
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.86'
$ws.Range("E31").Value = '  -0.59%  '
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.84'
$ws.Range("E32").Value = '  +12.98%  '
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.99'
$ws.Range("E33").Value = '  +3.27%  '
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0₃0814'
$ws.Range("E34").Value = '  +2.63%  '
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '173.01'
$ws.Range("E35").Value = '  -1.99%  '
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.11'
$ws.Range("E36").Value = '  +15.44%  '
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.407'
$ws.Range("E37").Value = '  +2.65%  '
# Row 38
$ws.Range("E38").Value = '  -0.09%  '
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '19.18'
$ws.Range("E39").Value = '  +1.90%  '
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.84'
$ws.Range("E40").Value = '  +8.99%  '
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '174.63'
$ws.Range("E41").Value = '  +10.95%  '
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.998'
$ws.Range("E42").Value = '  -0.12%  '
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '40.21'
$ws.Range("E43").Value = '  +0.14%  '
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.76'
$ws.Range("E44").Value = '  +2.20%  '
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '22.33'
$ws.Range("E45").Value = '  +6.64%  '
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0562'
$ws.Range("E46").Value = '  +5.66%  '
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.635'
$ws.Range("E47").Value = '  +0.82%  '
# Row 48
$ws.Range("B48").Value = 'VeChain'
$ws.Range("C48").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0241'
$ws.Range("E48").Value = '  +3.02%  '
# Row 49
$ws.Range("B49").Value = 'Stellar'
$ws.Range("C49").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0964'
$ws.Range("E49").Value = '  +0.86%  '
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '18.78'
$ws.Range("E50").Value = '  +4.89%  '
# Row 51
$ws.Range("B51").Value = 'dogwifhat'
$ws.Range("C51").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.72'
$ws.Range("E51").Value = '  +3.30%  '
